$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New words to append (in the order they were typed, before sorting)
$newWords = @(
    "xanadu",
    "cowabunga",
    "wannabe",
    "radical",
    "cool",
    "duh",
    "cd",
    "nerd",
    "dude",
    "chill",
    "whatever",
    "dynamite",
    "gnarly ",
    "disco",
    "groove"
)

$startRow = 12
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i]
}

# Sort the word list (A2:A26) ascending, keeping header in row 1
$sortRange = $ws.Range("A2:A26")
$keyRange = $ws.Range("A26")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Update the view: scroll and selection to match the final state
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B20").Select() | Out-Null
